$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 11.0
$ws.Columns("D").ColumnWidth = 8.333333333333332
$ws.Columns("E").ColumnWidth = 81.33333333333334
$ws.Columns("F").ColumnWidth = 5.666666666666667
$ws.Columns("G").ColumnWidth = 13.833333333333332

# ---------------------------------------------------------------------------
# Header row (row 2) : B2:E2 + G2
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "No. Día"
$ws.Range("C2").Value = "Fecha"
$ws.Range("D2").Value = "H-total"
$ws.Range("E2").Value = "Descipción"
$ws.Range("G2").Value = "Horas Total"
$ws.Range("G3").Value = "Restan"
$ws.Range("G4").Value = "Dias de trabajo"

$hdr = $ws.Range("B2:E2")
$hdr.Font.ThemeColor = 2
$hdr.Interior.ThemeColor = 4
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

$hdr2 = $ws.Range("G2:G4")
$hdr2.Font.ThemeColor = 2
$hdr2.Interior.ThemeColor = 4
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# H column values / formulas
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = 500
$ws.Range("H3").Formula = "=H2 - (SUM(D3:D23))"
$ws.Range("H4").Formula = "=SUM(B3:B23)"

# ---------------------------------------------------------------------------
# Data rows 3 .. 23 (B,C,D,E) with alternating banding
# row 3 = first data row (odd banding), row 4 = second data row (even banding), ...
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 1
$ws.Range("C3").Formula = "=TODAY()"
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = "Curso de Git, Github/ Instalacion de SQL SERVER / Conexión a la base de datos ACC MEX"

for ($r = 3; $r -le 23; $r++) {
    $bcd = $ws.Range("B" + $r + ":D" + $r)
    $e = $ws.Range("E" + $r)
    if (($r % 2) -eq 1) {
        # odd data row -> accent3 banding
        $bcd.Interior.ThemeColor = 7
        $bcd.HorizontalAlignment = -4108
        $bcd.VerticalAlignment = -4108
        $e.Interior.ThemeColor = 7
    } else {
        # even data row -> light1 banding
        $bcd.Interior.ThemeColor = 2
        $bcd.HorizontalAlignment = -4108
        $bcd.VerticalAlignment = -4108
        $e.Interior.ThemeColor = 2
    }
}

# Date formatting for C3 / C4 (banded like the rest of their row)
$ws.Range("C3").NumberFormat = "m/d/yyyy"
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C4").NumberFormat = "m/d/yyyy"
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("H5").Select()
